$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) values for rows 2-5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1033
$wsExhibit.Range("F3").Value = 2185
$wsExhibit.Range("F4").Value = 14
$wsExhibit.Range("F5").Value = 478

# Sheet "全部类型" - update column F (想去人数) values for rows 4-7
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1033
$wsAll.Range("F5").Value = 2185
$wsAll.Range("F6").Value = 14
$wsAll.Range("F7").Value = 478
